$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '29.059.47'
$ws.Range("E2").Value = '  +0.07%  '

$ws.Range("D3").Value = '1.833.76'
$ws.Range("E3").Value = '  +0.27%  '

$ws.Range("D4").Value = '''0.9997'
$ws.Range("E4").Value = '  +0.03%  '

$ws.Range("E5").Value = '  -0.83%  '

$ws.Range("D6").Value = '''0.6182'
$ws.Range("E6").Value = '  -2.08%  '

$ws.Range("E7").Value = '  +0.16%  '

$ws.Range("D8").Value = '''0.07461'
$ws.Range("E8").Value = '  -1.01%  '

$ws.Range("D9").Value = '''0.2925'
$ws.Range("E9").Value = '  -0.53%  '

$ws.Range("D10").Value = '''23.07'
$ws.Range("E10").Value = '  -0.20%  '

$ws.Range("D11").Value = '''0.07679'
$ws.Range("E11").Value = '  -0.33%  '

$ws.Range("D12").Value = '1.839.71'
$ws.Range("E12").Value = '  +0.62%  '

$ws.Range("D13").Value = '''5.001'
$ws.Range("E13").Value = '  +0.03%  '

$ws.Range("D14").Value = '''0.6717'
$ws.Range("E14").Value = '  +0.31%  '

$ws.Range("D15").Value = '''82.61'
$ws.Range("E15").Value = '  -0.65%  '

$ws.Range("D16").Value = '''0.000009184'
$ws.Range("E16").Value = '  -4.01%  '

$ws.Range("D17").Value = '''5.905'
$ws.Range("E17").Value = '  -2.86%  '

$ws.Range("D18").Value = '29.048.52'
$ws.Range("E18").Value = '  -0.05%  '

$ws.Range("D19").Value = '2.077.62'
$ws.Range("E19").Value = '  +0.00%  '

$ws.Range("D20").Value = '''232.48'
$ws.Range("E20").Value = '  +2.52%  '

$ws.Range("D21").Value = '''12.68'
$ws.Range("E21").Value = '  +0.86%  '

$ws.Range("E22").Value = '  +0.29%  '

$ws.Range("D23").Value = '''7.202'
$ws.Range("E23").Value = '  +0.83%  '

$ws.Range("D25").Value = '''159.31'
$ws.Range("E25").Value = '  -0.52%  '

$ws.Range("D26").Value = '''0.1417'
$ws.Range("E26").Value = '  -0.74%  '

$ws.Range("D27").Value = '''8.477'
$ws.Range("E27").Value = '  -0.39%  '

$ws.Range("D28").Value = '''17.79'

$ws.Range("D29").Value = '''1.501'
$ws.Range("E29").Value = '  -0.39%  '

$ws.Range("D30").Value = '''4.155'
$ws.Range("E30").Value = '  +0.21%  '

$ws.Range("D31").Value = '''4.104'
$ws.Range("E31").Value = '  +0.79%  '

$ws.Range("D32").Value = '''0.05536'
$ws.Range("E32").Value = '  +0.79%  '

$ws.Range("D34").Value = '''1.835'
$ws.Range("E34").Value = '  -1.22%  '

$ws.Range("D35").Value = '''0.7381'
$ws.Range("E35").Value = '  -0.89%  '

$ws.Range("D36").Value = '''1.139'
$ws.Range("E36").Value = '  -0.05%  '

$ws.Range("D37").Value = '''2.662'
$ws.Range("E37").Value = '  +0.24%  '

$ws.Range("D38").Value = '''2.774'
$ws.Range("E38").Value = '  +0.57%  '

$ws.Range("D39").Value = '''0.01779'
$ws.Range("E39").Value = '  -0.30%  '

$ws.Range("D40").Value = '1.208.31'
$ws.Range("E40").Value = '  -3.05%  '

$ws.Range("D41").Value = '''6.457'
$ws.Range("E41").Value = '  -2.07%  '

$ws.Range("D42").Value = '''0.8921'
$ws.Range("E42").Value = '  -1.23%  '

$ws.Range("E43").Value = '  +0.11%  '

$ws.Range("D44").Value = '''101.84'
$ws.Range("E44").Value = '  +0.39%  '

$ws.Range("D45").Value = '1.979.19'
$ws.Range("E45").Value = '  -0.07%  '

$ws.Range("D46").Value = '''65.43'
$ws.Range("E46").Value = '  +0.62%  '

$ws.Range("E47").Value = '  +0.22%  '

$ws.Range("D48").Value = '''0.5088'
$ws.Range("E48").Value = '  -0.27%  '

$ws.Range("D49").Value = '''0.4070'
$ws.Range("E49").Value = '  +0.10%  '

$ws.Range("D50").Value = '''9.152'
$ws.Range("E50").Value = '  +1.87%  '

$ws.Range("D51").Value = '''0.05816'
$ws.Range("E51").Value = '  +0.47%  '

